# Regenerate the "K" column (column G) of the save_data sheet.
# Previously column G held a different "Strike#" derived value; this
# recomputes/writes the actual strikeout (K) counts per outing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 2
    20 = 1
    21 = 0
    22 = 3
    23 = 1
    24 = 1
    25 = 2
    26 = 2
    27 = 0
    28 = 2
    29 = 2
    30 = 0
    31 = 1
    32 = 2
    33 = 0
    34 = 1
    35 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 3
    41 = 0
    42 = 1
    43 = 1
    44 = 0
    45 = 3
    46 = 1
    47 = 1
    48 = 0
    49 = 1
    50 = 2
    51 = 2
    52 = 1
    53 = 2
    54 = 0
    55 = 0
    56 = 0
    57 = 2
    58 = 1
    59 = 2
    60 = 2
    61 = 1
    62 = 3
    63 = 1
    64 = 1
    65 = 1
    66 = 1
    67 = 0
    68 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
